$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style (bold, border, alignment) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J60
$iVals = @(7,9,7,6,5,4,6,8,8,7,7,7,7,7,5,8,8,5,5,9,9,7,8,9,8,9,9,8,8,9,8,9,8,6,8,8,8,7,7,7,8,8,8,8,8,6,6,7,9,6,11,6,7,9,5,6,6,5,3)
$jVals = @(7,9,8,6,6,4,6,8,8,7,8,7,7,8,6,8,8,6,6,9,9,8,9,9,8,9,9,8,8,9,9,9,8,7,8,8,8,8,7,8,8,8,8,8,9,7,6,7,9,6,11,6,7,9,6,7,7,5,3)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}

$wb.Save()